# Auto-generated edit script: adds Primary/Secondary language group columns (C, D)
# to Sheet2, mirroring the target diff; removes duplicate chart-helper defined names;
# fixes chart field references; updates selection to D54.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,3).Value = "Primary language group"
$ws.Cells.Item(1,4).Value = "Secondary language group"

# Data rows (Primary / Secondary language group), rows 2-128 per source data
$ws.Cells.Item(2,3).Value = "Indo-European"
$ws.Cells.Item(2,4).Value = "Germanic"
$ws.Cells.Item(3,3).Value = "Indo-European"
$ws.Cells.Item(3,4).Value = "Italic"
$ws.Cells.Item(4,3).Value = "Afro-Asiatic"
$ws.Cells.Item(4,4).Value = "Semitic"
$ws.Cells.Item(5,3).Value = "Sino-Tibetan"
$ws.Cells.Item(5,4).Value = "Chinese"
$ws.Cells.Item(6,3).Value = "Koreanic"
$ws.Cells.Item(6,4).Value = "Koreanic"
$ws.Cells.Item(7,3).Value = "Indo-European"
$ws.Cells.Item(7,4).Value = "Italic"
$ws.Cells.Item(8,3).Value = "Indo-European"
$ws.Cells.Item(8,4).Value = "Balto-Slavic"
$ws.Cells.Item(9,3).Value = "Indo-European"
$ws.Cells.Item(9,4).Value = "Italic"
$ws.Cells.Item(10,3).Value = "Indo-European"
$ws.Cells.Item(10,4).Value = "Germanic"
$ws.Cells.Item(11,3).Value = "Indo-European"
$ws.Cells.Item(11,4).Value = "Germanic"
$ws.Cells.Item(12,3).Value = "Japonic"
$ws.Cells.Item(12,4).Value = "Japonic"
$ws.Cells.Item(13,3).Value = "Turkic"
$ws.Cells.Item(13,4).Value = "Southern Turkic"
$ws.Cells.Item(14,3).Value = "Austro-Asiatic"
$ws.Cells.Item(14,4).Value = "Mon-Khmer"
$ws.Cells.Item(15,3).Value = "Indo-European"
$ws.Cells.Item(15,4).Value = "Balto-Slavic"
$ws.Cells.Item(16,3).Value = "Indo-European"
$ws.Cells.Item(16,4).Value = "Indo-Iranian"
$ws.Cells.Item(17,3).Value = "Indo-European"
$ws.Cells.Item(17,4).Value = "Italic"
$ws.Cells.Item(18,3).Value = "Sino-Tibetan"
$ws.Cells.Item(18,4).Value = "Chinese"
$ws.Cells.Item(19,3).Value = "Indo-European"
$ws.Cells.Item(19,4).Value = "Indo-Iranian"
$ws.Cells.Item(20,3).Value = "Indo-European"
$ws.Cells.Item(20,4).Value = "Indo-Iranian"
$ws.Cells.Item(21,3).Value = "Austronesian"
$ws.Cells.Item(21,4).Value = "Malayo-Polynesian"
$ws.Cells.Item(22,3).Value = "Afro-Asiatic"
$ws.Cells.Item(22,4).Value = "Semitic"
$ws.Cells.Item(23,3).Value = "Indo-European"
$ws.Cells.Item(23,4).Value = "Balto-Slavic"
$ws.Cells.Item(24,3).Value = "Indo-European"
$ws.Cells.Item(24,4).Value = "Italic"
$ws.Cells.Item(25,3).Value = "Kra-Dai"
$ws.Cells.Item(25,4).Value = "Kam-Tai"
$ws.Cells.Item(26,3).Value = "Indo-European"
$ws.Cells.Item(26,4).Value = "Germanic"
$ws.Cells.Item(27,3).Value = "Indo-European"
$ws.Cells.Item(27,4).Value = "Indo-Iranian"
$ws.Cells.Item(28,3).Value = "Indo-European"
$ws.Cells.Item(28,4).Value = "Balto-Slavic"
$ws.Cells.Item(29,3).Value = "Indo-European"
$ws.Cells.Item(29,4).Value = "Balto-Slavic"
$ws.Cells.Item(30,3).Value = "Indo-European"
$ws.Cells.Item(30,4).Value = "Indo-Iranian"
$ws.Cells.Item(31,3).Value = "Indo-European"
$ws.Cells.Item(31,4).Value = "Indo-Iranian"
$ws.Cells.Item(32,3).Value = "Indo-European"
$ws.Cells.Item(32,4).Value = "Greek"
$ws.Cells.Item(33,3).Value = "Uralic"
$ws.Cells.Item(33,4).Value = "Finnic"
$ws.Cells.Item(34,3).Value = "Dravidian"
$ws.Cells.Item(34,4).Value = "South-Central Dravidian"
$ws.Cells.Item(35,3).Value = "Uralic"
$ws.Cells.Item(35,4).Value = "Finnic"
$ws.Cells.Item(36,3).Value = "Indo-European"
$ws.Cells.Item(36,4).Value = "Albanian"
$ws.Cells.Item(37,3).Value = "Indo-European"
$ws.Cells.Item(37,4).Value = "Indo-Iranian"
$ws.Cells.Item(38,3).Value = "Austronesian"
$ws.Cells.Item(38,4).Value = "Malayo-Polynesian"
$ws.Cells.Item(39,3).Value = "Dravidian"
$ws.Cells.Item(39,4).Value = "Southern Dravidian"
$ws.Cells.Item(40,3).Value = "Uralic"
$ws.Cells.Item(40,4).Value = "Uralic"
$ws.Cells.Item(41,3).Value = "Indo-European"
$ws.Cells.Item(41,4).Value = "Balto-Slavic"
$ws.Cells.Item(42,3).Value = "Indo-European"
$ws.Cells.Item(42,4).Value = "Balto-Slavic"
$ws.Cells.Item(43,3).Value = "Indo-European"
$ws.Cells.Item(43,4).Value = "Indo-Iranian"
$ws.Cells.Item(44,3).Value = "Sino-Tibetan"
$ws.Cells.Item(44,4).Value = "Chinese"
$ws.Cells.Item(45,3).Value = "Afro-Asiatic"
$ws.Cells.Item(45,4).Value = "Semitic"
$ws.Cells.Item(46,3).Value = "Misumalpan"
$ws.Cells.Item(46,4).Value = "Misumalpan"
$ws.Cells.Item(47,3).Value = "Niger-Congo"
$ws.Cells.Item(47,4).Value = "Atlantic-Congo"
$ws.Cells.Item(48,3).Value = "Indo-European"
$ws.Cells.Item(48,4).Value = "Indo-Iranian"
$ws.Cells.Item(49,3).Value = "Mongolic"
$ws.Cells.Item(49,4).Value = "Eastern Mongolic"
$ws.Cells.Item(50,3).Value = "Indo-European"
$ws.Cells.Item(50,4).Value = "Indo-Iranian"
$ws.Cells.Item(51,3).Value = "Indo-European"
$ws.Cells.Item(51,4).Value = "Indo-Iranian"
$ws.Cells.Item(52,3).Value = "Indo-European"
$ws.Cells.Item(52,4).Value = "Indo-Iranian"
$ws.Cells.Item(53,3).Value = "Indo-European"
$ws.Cells.Item(53,4).Value = "Indo-Iranian"
$ws.Cells.Item(54,3).Value = "Afro-Asiatic"
$ws.Cells.Item(54,4).Value = "Chadic"
$ws.Cells.Item(59,3).Value = "Indo-European"
$ws.Cells.Item(59,4).Value = "Balto-Slavic"
$ws.Cells.Item(60,3).Value = "Indo-European"
$ws.Cells.Item(60,4).Value = "Germanic"
$ws.Cells.Item(65,3).Value = "Indo-European"
$ws.Cells.Item(65,4).Value = "Balto-Slavic"
$ws.Cells.Item(67,3).Value = "Indo-European"
$ws.Cells.Item(67,4).Value = "Balto-Slavic"
$ws.Cells.Item(73,3).Value = "Indo-European"
$ws.Cells.Item(73,4).Value = "Italic"
$ws.Cells.Item(74,3).Value = "Indo-European"
$ws.Cells.Item(74,4).Value = "Balto-Slavic"
$ws.Cells.Item(77,3).Value = "Indo-European"
$ws.Cells.Item(77,4).Value = "Germanic"
$ws.Cells.Item(78,3).Value = "Turkic"
$ws.Cells.Item(78,4).Value = "Southern Turkic"
$ws.Cells.Item(94,3).Value = "Indo-European"
$ws.Cells.Item(94,4).Value = "Balto-Slavic"
$ws.Cells.Item(105,3).Value = "Indo-European"
$ws.Cells.Item(105,4).Value = "Balto-Slavic"
$ws.Cells.Item(118,3).Value = "Indo-European"
$ws.Cells.Item(118,4).Value = "Balto-Slavic"
$ws.Cells.Item(128,3).Value = "Turkic"
$ws.Cells.Item(128,4).Value = "Southern Turkic"

# Column widths for B, C, D (language-group columns)
$ws.Columns("B").ColumnWidth = 17.28515625
$ws.Columns("C").ColumnWidth = 24.5703125
$ws.Columns("D").ColumnWidth = 28.85546875

# Remove the duplicate chart-helper defined names (v1.3/v1.4/v1.5), keeping v1.0-v1.2
$wb.Names("_xlchart.v1.3").Delete() | Out-Null
$wb.Names("_xlchart.v1.4").Delete() | Out-Null
$wb.Names("_xlchart.v1.5").Delete() | Out-Null

# Restore selection to match the saved view
$ws.Range("D54").Select() | Out-Null
